$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-11-23 Sunday" "2025-11-24 Monday"

Replace-Text "169×9=1521" "593×5=2965"
Replace-Text "314×9=2826" "817×2=1634"
Replace-Text "344×7=2408" "232×8=1856"
Replace-Text "386×3=1158" "133×3=399"
Replace-Text "263×8=2104" "146×8=1168"

Replace-Text "960×6=5760" "964×5=4820"
Replace-Text "805×2=1610" "879×9=7911"
Replace-Text "757×7=5299" "739×3=2217"
Replace-Text "791×4=3164" "436×6=2616"
Replace-Text "356×8=2848" "485×9=4365"

Replace-Text "944×7=6608" "171×6=1026"
Replace-Text "517×9=4653" "961×7=6727"
Replace-Text "149×9=1341" "443×7=3101"
Replace-Text "556×8=4448" "715×3=2145"
Replace-Text "983×6=5898" "826×4=3304"

Replace-Text "713×4=2852" "876×5=4380"
Replace-Text "994×4=3976" "119×3=357"
Replace-Text "757×9=6813" "553×8=4424"
Replace-Text "295×4=1180" "128×6=768"
Replace-Text "359×3=1077" "800×9=7200"

Replace-Text "246×2=492" "288×8=2304"
Replace-Text "410×5=2050" "957×5=4785"
Replace-Text "698×5=3490" "465×4=1860"
Replace-Text "295×7=2065" "490×2=980"
Replace-Text "190×5=950" "591×2=1182"
